$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential notice date from 2021-03-22 to 2021-03-23
$ws.Range("A80").Value = $ws.Range("A80").Value2 -replace [regex]::Escape("2021-03-22"), "2021-03-23"

# Update Weight (D) and Percent Change (E) values for rows 2-77
$ws.Range("D2").Value = 0.07624656742143043 ; $ws.Range("E2").Value = -0.006888726801199452
$ws.Range("D3").Value = 0.04613530934973116 ; $ws.Range("E3").Value = 0.00856030628087967
$ws.Range("D4").Value = 0.03616476218943412 ; $ws.Range("E4").Value = 0.006737573626001137
$ws.Range("D5").Value = 0.03290361635510437 ; $ws.Range("E5").Value = 0.0007211538461537881
$ws.Range("D6").Value = 0.0313451436850157 ; $ws.Range("E6").Value = -0.01000198714976475
$ws.Range("D7").Value = 0.03011585548203736 ; $ws.Range("E7").Value = 0.005239598363117937
$ws.Range("D8").Value = 0.03015011360423401 ; $ws.Range("E8").Value = -0.0009345794392523477
$ws.Range("D9").Value = 0.02846578064682548 ; $ws.Range("E9").Value = -0.008188208979070244
$ws.Range("D10").Value = 0.02651588544924324 ; $ws.Range("E10").Value = -0.007733265720081262
$ws.Range("D11").Value = 0.02702862028823835 ; $ws.Range("E11").Value = 0.01858964510677508
$ws.Range("D12").Value = 0.02478824785252121 ; $ws.Range("E12").Value = -0.02141449756299929
$ws.Range("D13").Value = 0.0234574705804668 ; $ws.Range("E13").Value = -0.0201805629314924
$ws.Range("D14").Value = 0.02028792691126027 ; $ws.Range("E14").Value = 0.0003216374269006117
$ws.Range("D15").Value = 0.01892353416576615 ; $ws.Range("E15").Value = -0.02503657262277958
$ws.Range("D16").Value = 0.020638417656389 ; $ws.Range("E16").Value = -0.00576540755467192
$ws.Range("D17").Value = 0.01832221266735346 ; $ws.Range("E17").Value = 0.01186069351061403
$ws.Range("D18").Value = 0.01771139974114602 ; $ws.Range("E18").Value = 0.01204644412191569
$ws.Range("D19").Value = 0.01520136200011474 ; $ws.Range("E19").Value = -0.01234126274369518
$ws.Range("D20").Value = 0.01413481223620041 ; $ws.Range("E20").Value = -0.01000944287063266
$ws.Range("D21").Value = 0.01596210982476723 ; $ws.Range("E21").Value = -0.0099134700551885
$ws.Range("D22").Value = 0.01417386550681419 ; $ws.Range("E22").Value = -0.008733258928571463
$ws.Range("D23").Value = 0.01404316063655743 ; $ws.Range("E23").Value = -0.03928526169053348
$ws.Range("D24").Value = 0.01512192072685354 ; $ws.Range("E24").Value = 0.0003334444814937498
$ws.Range("D25").Value = 0.01450843833657883 ; $ws.Range("E25").Value = 0.004470370168456528
$ws.Range("D26").Value = 0.01263467229943359 ; $ws.Range("E26").Value = -0.06441717791411039
$ws.Range("D27").Value = 0.01227063649588943 ; $ws.Range("E27").Value = -0.02723390540649417
$ws.Range("D28").Value = 0.01265311137530566 ; $ws.Range("E28").Value = -0.03275940880694805
$ws.Range("D29").Value = 0.01226114506809469 ; $ws.Range("E29").Value = -0.0040801844952989
$ws.Range("D30").Value = 0.01107392564143576 ; $ws.Range("E30").Value = -0.02789136296269845
$ws.Range("D31").Value = 0.01215199364845516 ; $ws.Range("E31").Value = -0.02684891383939469
$ws.Range("D32").Value = 0.0128552491266218 ; $ws.Range("E32").Value = 0.003599363189589644
$ws.Range("D33").Value = 0.01123794937801364 ; $ws.Range("E33").Value = -0.005828531210135046
$ws.Range("D34").Value = 0.01146974784118834 ; $ws.Range("E34").Value = 0.005654709312599371
$ws.Range("D35").Value = 0.009936338472620159 ; $ws.Range("E35").Value = -0.01170149253731345
$ws.Range("D36").Value = 0.01084099018430647 ; $ws.Range("E36").Value = 0.00764705882352934
$ws.Range("D37").Value = 0.01083159762555126 ; $ws.Range("E37").Value = -0.01086212404728248
$ws.Range("D38").Value = 0.01000010900311608 ; $ws.Range("E38").Value = -0.01431113747590096
$ws.Range("D39").Value = 0.009398095421426684 ; $ws.Range("E39").Value = -0.02360713684564886
$ws.Range("D40").Value = 0.009607697785227229 ; $ws.Range("E40").Value = -0.008067835679591684
$ws.Range("D41").Value = 0.009280144657269115 ; $ws.Range("E41").Value = -0.01170855395629788
$ws.Range("D42").Value = 0.009569435466929678 ; $ws.Range("E42").Value = -0.008317060823027522
$ws.Range("D43").Value = 0.009962341030016171 ; $ws.Range("E43").Value = -0.0159979357502259
$ws.Range("D44").Value = 0.009179792582147626 ; $ws.Range("E44").Value = 0.01432448733413749
$ws.Range("D45").Value = 0.009124228181932576 ; $ws.Range("E45").Value = -0.005558806319485021
$ws.Range("D46").Value = 0.009610070642175915 ; $ws.Range("E46").Value = -0.01777777777777778
$ws.Range("D47").Value = 0.0088317735630071 ; $ws.Range("E47").Value = -0.02364320257925834
$ws.Range("D48").Value = 0.007349528922394947 ; $ws.Range("E48").Value = -0.01162290142057687
$ws.Range("D49").Value = 0.00829323390469504 ; $ws.Range("E49").Value = -0.01096791883740056
$ws.Range("D50").Value = 0.007976506936565302 ; $ws.Range("E50").Value = -0.01227107929720184
$ws.Range("D51").Value = 0.007821233109985701 ; $ws.Range("E51").Value = -0.03437748871781265
$ws.Range("D52").Value = 0.007627103750871375 ; $ws.Range("E52").Value = -0.02136926636722458
$ws.Range("D53").Value = 0.007114121739277448 ; $ws.Range("E53").Value = -0.04857202418178019
$ws.Range("D54").Value = 0.007226535837221419 ; $ws.Range("E54").Value = -0.01963279155037478
$ws.Range("D55").Value = 0.006691456595292859 ; $ws.Range("E55").Value = -0.01470855496453904
$ws.Range("D56").Value = 0.006369539002587872 ; $ws.Range("E56").Value = -0.007450639513224777
$ws.Range("D57").Value = 0.006897895149828491 ; $ws.Range("E57").Value = -0.01719986240110083
$ws.Range("D58").Value = 0.006403154476027582 ; $ws.Range("E58").Value = -0.01459143968871601
$ws.Range("D59").Value = 0.005764065004514978 ; $ws.Range("E59").Value = -0.05415094339622628
$ws.Range("D60").Value = 0.006265083862325949 ; $ws.Range("E60").Value = -0.0124275062137531
$ws.Range("D61").Value = 0.005548036153145077 ; $ws.Range("E61").Value = -0.05016484006058985
$ws.Range("D62").Value = 0.005698811438426129 ; $ws.Range("E62").Value = 0.01651630811936156
$ws.Range("D63").Value = 0.005261019331393671 ; $ws.Range("E63").Value = 0.01691347816282041
$ws.Range("D64").Value = 0.005013055780256045 ; $ws.Range("E64").Value = -0.01199116440517523
$ws.Range("D65").Value = 0.004777747466178075 ; $ws.Range("E65").Value = -0.005628673123085903
$ws.Range("D66").Value = 0.004330463931350875 ; $ws.Range("E66").Value = -0.008835616438356175
$ws.Range("D67").Value = 0.004476987847932199 ; $ws.Range("E67").Value = -0.02729561415131843
$ws.Range("D68").Value = 0.003624316533778026 ; $ws.Range("E68").Value = -0.002516521063076538
$ws.Range("D69").Value = 0.004063814131742352 ; $ws.Range("E69").Value = -0.0263119480329902
$ws.Range("D70").Value = 0.003819508735067284 ; $ws.Range("E70").Value = -0.04014806378132119
$ws.Range("D71").Value = 0.003191047685304 ; $ws.Range("E71").Value = -0.0007900729655621053
$ws.Range("D72").Value = 0.002607423744966917 ; $ws.Range("E72").Value = -0.008759124087591275
$ws.Range("D73").Value = 0.00258596916338922 ; $ws.Range("E73").Value = 0.02290149299382538
$ws.Range("D74").Value = 0.002236467108655764 ; $ws.Range("E74").Value = 0.0172188943657301
$ws.Range("D75").Value = 0.001936449008206353 ; $ws.Range("E75").Value = -0.0521801286633311
$ws.Range("D76").Value = 0.001901844844371357 ; $ws.Range("E76").Value = -0.0782387190684134
$ws.Range("E77").Value = -0.008408339712240198

$ws.Protect()
Write-Host "Edit complete"